# "sumproduct() working on example.xlsx"
#
# Sheet2 used to hold a single SUM() formula pulling from Sheet1. Replace it
# with a small 3x3 sample table (columns A/B/C, rows 1-3) and a SUMPRODUCT()
# formula in E2 that multiplies the A/B columns together, gated by a C<9
# condition.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 1
$ws2.Range("A1").Value = 3
$ws2.Range("B1").Value = 4
$ws2.Range("C1").Value = 7

# Row 2 (+ the new SUMPRODUCT formula)
$ws2.Range("A2").Value = 2
$ws2.Range("B2").Value = 5
$ws2.Range("C2").Value = 8
$ws2.Range("E2").Formula = "=SUMPRODUCT(A1:A3*B1:B3*(C1:C3<9))"

# Row 3
$ws2.Range("A3").Value = 1
$ws2.Range("B3").Value = 6
$ws2.Range("C3").Value = 9

# Cursor ends up on Sheet1!A2, but Sheet2 remains the active tab.
$ws1.Range("A2").Select() | Out-Null
$ws2.Select() | Out-Null
